# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns
# per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.504.92"
$ws.Range("E2").Value = "'  +0.72%  "
$ws.Range("D3").Value = "'1.646.58"
$ws.Range("E3").Value = "'  +0.87%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E5").Value = "'  +0.26%  "
$ws.Range("D6").Value = "'302.45"
$ws.Range("E6").Value = "'  -0.17%  "
$ws.Range("D7").Value = "'0.3842"
$ws.Range("E7").Value = "'  +0.70%  "
$ws.Range("D8").Value = "'0.3598"
$ws.Range("E8").Value = "'  +1.30%  "
$ws.Range("D9").Value = "'51.02"
$ws.Range("E9").Value = "'  -1.87%  "
$ws.Range("D10").Value = "'0.08170"
$ws.Range("E10").Value = "'  +0.51%  "
$ws.Range("D11").Value = "'1.231"
$ws.Range("E11").Value = "'  +0.70%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "'  +0.09%  "
$ws.Range("D13").Value = "'22.32"
$ws.Range("E13").Value = "'  +0.32%  "
$ws.Range("D14").Value = "'6.450"
$ws.Range("E14").Value = "'  +0.38%  "
$ws.Range("D15").Value = "'7.455"
$ws.Range("E15").Value = "'  +2.26%  "
$ws.Range("D16").Value = "'0.00001222"
$ws.Range("E16").Value = "'  -0.24%  "
$ws.Range("D17").Value = "'1.647.17"
$ws.Range("E17").Value = "'  +1.26%  "
$ws.Range("D18").Value = "'97.54"
$ws.Range("E18").Value = "'  +3.01%  "
$ws.Range("D19").Value = "'0.07015"
$ws.Range("E19").Value = "'  +1.21%  "
$ws.Range("D20").Value = "'6.770"
$ws.Range("E20").Value = "'  +2.91%  "
$ws.Range("D21").Value = "'17.56"
$ws.Range("E21").Value = "'  +1.42%  "
$ws.Range("E22").Value = "'  +0.18%  "
$ws.Range("D23").Value = "'12.61"
$ws.Range("E23").Value = "'  +1.72%  "
$ws.Range("D24").Value = "'23.516.92"
$ws.Range("E24").Value = "'  +0.79%  "
$ws.Range("D25").Value = "'2.482"
$ws.Range("E25").Value = "'  -2.83%  "
$ws.Range("D26").Value = "'3.038"
$ws.Range("E26").Value = "'  -2.84%  "
$ws.Range("D27").Value = "'21.22"
$ws.Range("E27").Value = "'  +1.04%  "
$ws.Range("D28").Value = "'153.65"
$ws.Range("E28").Value = "'  +1.37%  "
$ws.Range("D29").Value = "'5.232"
$ws.Range("E29").Value = "'  -0.68%  "
$ws.Range("D30").Value = "'133.98"
$ws.Range("E30").Value = "'  +0.77%  "
$ws.Range("D31").Value = "'1.826.58"
$ws.Range("E31").Value = "'  +0.99%  "
$ws.Range("D32").Value = "'7.128"
$ws.Range("E32").Value = "'  +9.68%  "
$ws.Range("D33").Value = "'2.243"
$ws.Range("E33").Value = "'  +4.53%  "
$ws.Range("D34").Value = "'12.21"
$ws.Range("E34").Value = "'  +5.69%  "
$ws.Range("D35").Value = "'1.055"
$ws.Range("E35").Value = "'  -1.74%  "
$ws.Range("D36").Value = "'0.02792"
$ws.Range("E36").Value = "'  +2.04%  "
$ws.Range("D37").Value = "'0.2503"
$ws.Range("E37").Value = "'  +0.64%  "
$ws.Range("D38").Value = "'0.08772"
$ws.Range("D39").Value = "'6.075"
$ws.Range("E39").Value = "'  +2.55%  "
$ws.Range("D40").Value = "'0.06977"
$ws.Range("E40").Value = "'  +0.26%  "
$ws.Range("D41").Value = "'13.08"
$ws.Range("E41").Value = "'  +7.97%  "
$ws.Range("D42").Value = "'0.6984"
$ws.Range("E42").Value = "'  +0.62%  "
$ws.Range("D43").Value = "'1.335"
$ws.Range("E43").Value = "'  +1.29%  "
$ws.Range("D44").Value = "'15.95"
$ws.Range("E44").Value = "'  +3.76%  "
$ws.Range("D45").Value = "'0.6506"
$ws.Range("E45").Value = "'  +1.90%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "'  +0.23%  "
$ws.Range("D47").Value = "'2.295"
$ws.Range("E47").Value = "'  +1.22%  "
$ws.Range("D48").Value = "'3.955"
$ws.Range("E48").Value = "'  -0.03%  "
$ws.Range("D49").Value = "'0.07876"
$ws.Range("E49").Value = "'  -0.60%  "
$ws.Range("D50").Value = "'128.05"
$ws.Range("E50").Value = "'  -1.18%  "
$ws.Range("D51").Value = "'1.177"
$ws.Range("E51").Value = "'  -0.37%  "
